{"js": "// Fix a typo in the \"Traversing the B+ tree\" section: \"grater\" -> \"greater\".\n// (Word re-chunked the surrounding runs during its spell-check pass when the\n// file was resaved, but the only actual content change is this one word.)\nconst body = context.document.body;\n\nconst results = body.search(\"grater\", {\n  matchCase: true,\n  matchWholeWord: true,\n});\nresults.load(\"items/text\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(\"greater\", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Fix a typo in the \"Traversing the B+ tree\" section: \"grater\" -> \"greater\".\n# (Word re-chunked the surrounding runs during its spell-check pass when the\n# file was resaved, but the only actual content change is this one word.)\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"grater\"\n$find.Replacement.Text = \"greater\"\n$find.MatchCase = $true\n$find.MatchWholeWord = $true\n$find.Wrap = 1\n$find.Execute($null, $true, $true, $false, $null, $null, $true, 1, $false, \"greater\", 2)\n"}
